# UK removed in RP3
# The "ERT_SU_CZ" sheet has a per-country data table (rows 7-36) summarised
# in row 6. The United Kingdom occupied row 36 - delete that entire row so
# the sheet, the summary formulas (which auto-adjust their SUM/shared
# ranges from B7:B36 to B7:B35, etc.) and the shared-strings table (the now
# unused "United Kingdom" string is dropped, shifting later entries such as
# the "Change Log" sheet headers down by one index) all update together,
# exactly as Excel would do it interactively.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ERT_SU_CZ")

$ws.Rows.Item(36).Delete()
